# The SerializationPath cell (D8 on "General Settings") currently holds a
# formula that calls into the external add-in workbook FixedIncome.xla
# (`=[1]!qlSerializationPath(Trigger)`). Per the commit message we need to
# remove the link to that external workbook: replace the formula with a
# plain literal text value and drop the external reference/link entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the formula with a literal text value. The leading apostrophe is
# Excel's "treat as text" prefix marker (it is not stored as part of the
# cell's text) - it keeps the cell's existing text-quote-prefixed style
# instead of Excel re-deriving a fresh (unprefixed) style for the literal.
$ws.Range("D8").Value = "'C:\Users\erik\junk\"

# Now that nothing in the workbook references the external workbook anymore,
# break the link so the externalReference/externalLink parts are dropped
# from the saved package.
$linkSources = $wb.LinkSources(1)
if ($linkSources) {
    foreach ($link in $linkSources) {
        $wb.BreakLink($link, 1)
    }
}
